$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: move the "FORMATION ACADEMIQUE" section (header + its 4
# entries) so that it now appears AFTER the whole "EXPERIENCE
# PROFESSIONNELLE" section instead of before it (right before the
# "COMPETENCES TECHNIQUES" header).
# ------------------------------------------------------------------

# Locate the start paragraph of the "FORMATION ACADEMIQUE" block.
$formationStart = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "FORMATION ACADEMIQUE*") {
        $formationStart = $i
        break
    }
}

# Locate the "EXPERIENCE PROFESSIONNELLE" header paragraph; the
# formation block ends right before it.
$expIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "EXPERIENCE PROFESSIONNELLE*") {
        $expIdx = $i
        break
    }
}
$formationEnd = $expIdx - 1

# Locate the "COMPETENCES TECHNIQUES" header paragraph; this is where
# the formation block needs to be re-inserted, right before it.
$compIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "COMPETENCES TECHNIQUES*") {
        $compIdx = $i
        break
    }
}

# Insert each paragraph of the FORMATION ACADEMIQUE block, in its
# original order, right before "COMPETENCES TECHNIQUES". Paragraphs
# are copied one at a time (rather than as one multi-paragraph range)
# so each paragraph keeps its own paragraph formatting (pPr).
$insertPos = $d.Paragraphs.Item($compIdx).Range.Start
for ($i = $formationStart; $i -le $formationEnd; $i++) {
    $src = $d.Paragraphs.Item($i).Range
    $ft = $src.FormattedText
    $len = $src.End - $src.Start
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.FormattedText = $ft
    $insertPos = $insertPos + $len
}

# Delete the original FORMATION ACADEMIQUE block (its paragraph
# indices are unaffected by the insertion above, which happened
# further down in the document).
$delStart = $d.Paragraphs.Item($formationStart).Range.Start
$delEnd = $d.Paragraphs.Item($formationEnd).Range.End
$d.Range($delStart, $delEnd).Delete()

# ------------------------------------------------------------------
# Part 2: within "COMPETENCES TECHNIQUES", reorder the four lines
#   MLOps / Bases de données / Autres / Langages
# into
#   Autres / Langages / Bases de données / MLOps
# All four paragraphs share identical paragraph formatting, so the
# reorder is done by swapping the run text in place (leaving the
# paragraph marks/formatting untouched).
# ------------------------------------------------------------------

$mlopsIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "MLOps :*") {
        $mlopsIdx = $i
        break
    }
}

$bddIdx = $mlopsIdx + 1
$autresIdx = $mlopsIdx + 2
$langagesIdx = $mlopsIdx + 3

$mlopsText = "MLOps : aws, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"
$bddText = "Bases de données : SQL, MongoDB, Neo4j, Redis"
$autresText = "Autres : scikit-learn/pandas, postgrsql, au data cleaning"
$langagesText = "Langages : r, python, matlab, c, c++"

function Set-ParaText($idx, $newText) {
    $r = $d.Paragraphs.Item($idx).Range
    $noMark = $d.Range($r.Start, $r.End - 1)
    $noMark.Text = $newText
}

Set-ParaText $mlopsIdx $autresText
Set-ParaText $bddIdx $langagesText
Set-ParaText $autresIdx $bddText
Set-ParaText $langagesIdx $mlopsText
